# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.941.56"
$ws.Range("E2").Value = "  +0.20%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.796.91"
$ws.Range("E3").Value = "  -0.68%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "359.91"
$ws.Range("E5").Value = "  +1.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.17"
$ws.Range("E6").Value = "  -1.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.562"
$ws.Range("E7").Value = "  -0.64%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  -1.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.26"
$ws.Range("E10").Value = "  -1.36%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.134"
$ws.Range("E11").Value = "  +2.28%  "

$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0850"
$ws.Range("E12").Value = "  -0.50%  "

$ws.Range("E13").Value = "  -1.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.66"
$ws.Range("E14").Value = "  -1.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.233.25"
$ws.Range("E15").Value = "  -0.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.806.56"
$ws.Range("E16").Value = "  -1.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.945"
$ws.Range("E17").Value = "  +2.98%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.898.84"
$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.47"
$ws.Range("E19").Value = "  -1.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.08"
$ws.Range("E20").Value = "  -1.97%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.36"
$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("E22").Value = "  -1.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.35"
$ws.Range("E23").Value = "  +0.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.43"
$ws.Range("E24").Value = "  +1.10%  "

$ws.Range("E25").Value = "  -0.81%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.54"
$ws.Range("E26").Value = "  -1.55%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.166"
$ws.Range("E28").Value = "  +18.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.28"
$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("E30").Value = "  -3.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.15"

$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "35.05"
$ws.Range("E32").Value = "  +2.36%  "

$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "52.15"
$ws.Range("E33").Value = "  -0.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0468"
$ws.Range("E34").Value = "  -2.06%  "

$ws.Range("E35").Value = "  +0.68%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.18"
$ws.Range("E36").Value = "  -3.59%  "

$ws.Range("E37").Value = "  -0.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.87"
$ws.Range("E38").Value = "  +2.62%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.21"
$ws.Range("E39").Value = "  -2.85%  "

$ws.Range("E40").Value = "  -2.79%  "

$ws.Range("E41").Value = "  +3.13%  "

$ws.Range("E42").Value = "  -1.70%  "

$ws.Range("E43").Value = "  -1.70%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.88"
$ws.Range("E44").Value = "  -6.35%  "

$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "119.54"
$ws.Range("E45").Value = "  -3.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.097.20"
$ws.Range("E46").Value = "  +0.27%  "

$ws.Range("E47").Value = "  -1.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.28"
$ws.Range("E48").Value = "  +0.95%  "

$ws.Range("E49").Value = "  -3.60%  "

$ws.Range("E50").Value = "  -2.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.13"
$ws.Range("E51").Value = "  +29.42%  "
